$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 1656
$ws.Range("J3").Value = 1740
$ws.Range("J4").Value = 388
$ws.Range("J5").Value = 120
$ws.Range("J6").Value = 2270
$ws.Range("J7").Value = 6174

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("J6").Value = 21
$ws.Range("J7").Value = 69

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("J3").Value = 28
$ws.Range("J7").Value = 79

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("J2").Value = 45
$ws.Range("J3").Value = 93
$ws.Range("J6").Value = 70
$ws.Range("J7").Value = 226

$ws = $wb.Worksheets.Item('New City')
$ws.Range("J3").Value = 40
$ws.Range("J4").Value = 9
$ws.Range("J6").Value = 65
$ws.Range("J7").Value = 160

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("J2").Value = 51
$ws.Range("J5").Value = 17
$ws.Range("J7").Value = 175
$ws.Range("J8").Value = 377
$ws.Range("J9").Value = 39
$ws.Range("J15").Value = 81
$ws.Range("J19").Value = 212
$ws.Range("J20").Value = 127
$ws.Range("J29").Value = 344
$ws.Range("J33").Value = 264
$ws.Range("J34").Value = 39
$ws.Range("J36").Value = 94
$ws.Range("J42").Value = 238
$ws.Range("J51").Value = 82
$ws.Range("J52").Value = 143
$ws.Range("J54").Value = 121
$ws.Range("J55").Value = 70
$ws.Range("J57").Value = 30
$ws.Range("J61").Value = 12
$ws.Range("J63").Value = 30
$ws.Range("J64").Value = 42
$ws.Range("J65").Value = 160
$ws.Range("J67").Value = 226
$ws.Range("J71").Value = 29
$ws.Range("J74").Value = 10
$ws.Range("J77").Value = 45
$ws.Range("J78").Value = 78
$ws.Range("J79").Value = 186
$ws.Range("J85").Value = 282
$ws.Range("J86").Value = 34
$ws.Range("J87").Value = 25
$ws.Range("J89").Value = 69
$ws.Range("J90").Value = 69
$ws.Range("J91").Value = 71
$ws.Range("J92").Value = 23
$ws.Range("J95").Value = 92
$ws.Range("J98").Value = 42
$ws.Range("J99").Value = 79
$ws.Range("J101").Value = 6174

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("J3").Value = 27
$ws.Range("J7").Value = 92

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("J2").Value = 66
$ws.Range("J3").Value = 74
$ws.Range("J7").Value = 264

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("J6").Value = 60
$ws.Range("J7").Value = 121

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("J2").Value = 101
$ws.Range("J3").Value = 123
$ws.Range("J5").Value = 13
$ws.Range("J6").Value = 89
$ws.Range("J7").Value = 344

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("J3").Value = 60
$ws.Range("J7").Value = 212

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("J2").Value = 72
$ws.Range("J3").Value = 111
$ws.Range("J7").Value = 282

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J6").Value = 128
$ws.Range("J7").Value = 238

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("J6").Value = 20
$ws.Range("J7").Value = 78

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("J2").Value = 20
$ws.Range("J7").Value = 70

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("J2").Value = 20
$ws.Range("J6").Value = 13
$ws.Range("J7").Value = 71

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("J2").Value = 50
$ws.Range("J3").Value = 66
$ws.Range("J6").Value = 54
$ws.Range("J7").Value = 186

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("J2").Value = 9
$ws.Range("J7").Value = 42

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("J6").Value = 38
$ws.Range("J7").Value = 127

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("J6").Value = 42
$ws.Range("J7").Value = 94

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("J3").Value = 47
$ws.Range("J6").Value = 55
$ws.Range("J7").Value = 143

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range("J2").Value = 12
$ws.Range("J7").Value = 39

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("J3").Value = 21
$ws.Range("J7").Value = 81

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("J4").Value = 4
$ws.Range("J7").Value = 42

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("J2").Value = 8
$ws.Range("J6").Value = 18
$ws.Range("J7").Value = 39

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("J2").Value = 19
$ws.Range("J7").Value = 51

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range("J6").Value = 9
$ws.Range("J7").Value = 23

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("J2").Value = 123
$ws.Range("J6").Value = 99
$ws.Range("J7").Value = 377

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("J2").Value = 4
$ws.Range("J7").Value = 17

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("J2").Value = 6
$ws.Range("J7").Value = 34

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("J3").Value = 18
$ws.Range("J7").Value = 69

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("J2").Value = 20
$ws.Range("J3").Value = 24
$ws.Range("J6").Value = 22
$ws.Range("J7").Value = 82

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("J3").Value = 8
$ws.Range("J7").Value = 30

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range("J6").Value = 16
$ws.Range("J7").Value = 29

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("J4").Value = 4
$ws.Range("J7").Value = 45

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("J2").Value = 56
$ws.Range("J6").Value = 60
$ws.Range("J7").Value = 175

$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Range("J6").Value = 12
$ws.Range("J7").Value = 25

$ws = $wb.Worksheets.Item('Mount Greenwood')
$ws.Range("J3").Value = 1
$ws.Range("J7").Value = 12

$ws = $wb.Worksheets.Item('Printers Row')
$ws.Range("J3").Value = 5
$ws.Range("J7").Value = 10
